# JS-SPA-Self-Evaluation-Protocol — "added functionality for delete add"
#
# Fills in previously-blank "Numbers of Commits in GitHub" score cells,
# marks the Admin Home Screen as partially done ("Yes Half"), updates the
# commit-history comments with more detail, grows the comment row so the
# extra text fits, resets the current view/selection, and lets the Total
# Score formula recalculate on its own.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Commit-history comments (row 8 / row 9) -------------------------------
# Row 8 needs to be taller to fit the longer comment.
$ws.Rows.Item(8).RowHeight = 60

$ws.Range("E8").Value = "I created half the project during the lab so half the project is committed at once`n - after the lab was finished, overall 4 days of work (except the lab)"
$ws.Range("E9").Value = "the project is commited on small portions`n - over 30 commits"

# --- Basic Options scores that were left blank -----------------------------
$ws.Range("C21").Value = 5   # List User Ads
$ws.Range("C22").Value = 5   # Show Buttons
$ws.Range("C25").Value = 5   # Deactivate/Publish Again Ad
$ws.Range("C26").Value = 10  # Edit Inactive Ads
$ws.Range("C27").Value = 5   # Change and Edit Images
$ws.Range("C28").Value = 5   # Delete Ad
$ws.Range("C32").Value = 5   # Authorization Checks

# --- Admin Options: Admin Home Screen implemented only halfway -------------
$ws.Range("C34").Value = "Yes Half"

# --- Refresh the current selection / scroll position -----------------------
$ws.Application.Goto($ws.Range("D34"))
$ws.Range("D34").Select()
$excel.ActiveWindow.ScrollRow = 31

$wb.Application.Calculate()
